$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 139, shifting existing rows 139:248 down to 140:249.
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row 139 with the new weekly record.
$ws.Cells.Item(139, 1).Value = 4
$ws.Cells.Item(139, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(139, 3).Value = "Los Lagos"
$ws.Cells.Item(139, 4).Value = (Get-Date -Year 2022 -Month 12 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(139, 5).Value = 10
$ws.Cells.Item(139, 6).Value = "Fruta"
$ws.Cells.Item(139, 7).Value = 100103
$ws.Cells.Item(139, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(139, 9).Value = 100103002
$ws.Cells.Item(139, 10).Value = "Ciruela"
$ws.Cells.Item(139, 11).Value = "Angeleno"
$ws.Cells.Item(139, 12).Value = "Primera"
$ws.Cells.Item(139, 13).Value = 600
$ws.Cells.Item(139, 14).Value = 17000
$ws.Cells.Item(139, 15).Value = 18000
$ws.Cells.Item(139, 16).Value = 17500
$ws.Cells.Item(139, 17).Value = "`$/caja 14 kilos granel"
$ws.Cells.Item(139, 18).Value = "Región Metropolitana"
$ws.Cells.Item(139, 19).Value = 1250
$ws.Cells.Item(139, 20).Value = 14
